$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: replace the "g"/"g" placeholder row with real admin (BR) data,
# matching the shape of the other data rows above it.
$ws.Range("A10").Value = "admin"
$ws.Range("B10").Value = "Admin (BR)"
$ws.Range("C10").Value = 123456
$ws.Range("D10").Value = "município"
$ws.Range("E10").Value = 2602902
$ws.Range("F10").Value = 10000
$ws.Rows.Item(10).RowHeight = 13.8

# Row 11: drop the old "gesto"/"GESTO" row entirely, leaving only the
# pre-existing style on F11 with no content.
$ws.Range("A11:F11").ClearContents()

# Move the active selection to D11.
$ws.Range("D11").Select()

# Minor page-setup touch-up (first page number reset), matching the
# resave that produced the new pageSetup attributes.
$ws.PageSetup.FirstPageNumber = 0
